# Append the final 12 games of the 2023-24 NBA regular season (attendance
# 17832, played same day) to the bottom of the results table on Sheet1.
#
# Column layout (row 1 header): A=Away team, B=Away Pts, C=Home team,
# D=Home Pts, E=Overtime, F=Attend., G=Arena, H=Win, I=Loss
# New rows occupy 823-834 (table previously ended at row 822).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newGames = @(
    @{ Row = 823; Away = 'Orlando Magic';         AwayPts = 116; Home = 'Cleveland Cavaliers';  HomePts = 109; OT = 'No'; Attend = 17832; Arena = 'Rocket Mortgage Fieldhouse'; Win = 'Orlando Magic';         Loss = 'Cleveland Cavaliers' }
    @{ Row = 824; Away = 'Detroit Pistons';        AwayPts = 115; Home = 'Indiana Pacers';       HomePts = 129; OT = 'No'; Attend = 17832; Arena = 'Gainbridge Fieldhouse';       Win = 'Indiana Pacers';        Loss = 'Detroit Pistons' }
    @{ Row = 825; Away = 'New York Knicks';        AwayPts = 110; Home = 'Philadelphia 76ers';   HomePts = 96;  OT = 'No'; Attend = 17832; Arena = 'Wells Fargo Center';          Win = 'New York Knicks';       Loss = 'Philadelphia 76ers' }
    @{ Row = 826; Away = 'Brooklyn Nets';          AwayPts = 93;  Home = 'Toronto Raptors';      HomePts = 121; OT = 'No'; Attend = 17832; Arena = 'Scotiabank Arena';            Win = 'Toronto Raptors';       Loss = 'Brooklyn Nets' }
    @{ Row = 827; Away = 'Phoenix Suns';           AwayPts = 113; Home = 'Dallas Mavericks';     HomePts = 123; OT = 'No'; Attend = 17832; Arena = 'American Airlines Center';    Win = 'Dallas Mavericks';      Loss = 'Phoenix Suns' }
    @{ Row = 828; Away = 'Boston Celtics';         AwayPts = 129; Home = 'Chicago Bulls';        HomePts = 112; OT = 'No'; Attend = 17832; Arena = 'United Center';               Win = 'Boston Celtics';        Loss = 'Chicago Bulls' }
    @{ Row = 829; Away = 'Houston Rockets';        AwayPts = 105; Home = 'New Orleans Pelicans'; HomePts = 127; OT = 'No'; Attend = 17832; Arena = 'Smoothie King Center';        Win = 'New Orleans Pelicans';  Loss = 'Houston Rockets' }
    @{ Row = 830; Away = 'Los Angeles Clippers';   AwayPts = 107; Home = 'Oklahoma City Thunder'; HomePts = 129; OT = 'No'; Attend = 17832; Arena = 'Paycom Center';              Win = 'Oklahoma City Thunder'; Loss = 'Los Angeles Clippers' }
    @{ Row = 831; Away = 'Washington Wizards';     AwayPts = 110; Home = 'Denver Nuggets';       HomePts = 130; OT = 'No'; Attend = 17832; Arena = 'Ball Arena';                  Win = 'Denver Nuggets';        Loss = 'Washington Wizards' }
    @{ Row = 832; Away = 'Charlotte Hornets';      AwayPts = 115; Home = 'Utah Jazz';            HomePts = 107; OT = 'No'; Attend = 17832; Arena = 'Delta Center';                Win = 'Charlotte Hornets';     Loss = 'Utah Jazz' }
    @{ Row = 833; Away = 'Los Angeles Lakers';     AwayPts = 110; Home = 'Golden State Warriors'; HomePts = 128; OT = 'No'; Attend = 17832; Arena = 'Chase Center';               Win = 'Golden State Warriors'; Loss = 'Los Angeles Lakers' }
    @{ Row = 834; Away = 'San Antonio Spurs';      AwayPts = 122; Home = 'Sacramento Kings';     HomePts = 127; OT = 'No'; Attend = 17832; Arena = 'Golden 1 Center';             Win = 'Sacramento Kings';      Loss = 'San Antonio Spurs' }
)

foreach ($game in $newGames) {
    $r = $game.Row
    $ws.Cells.Item($r, 1).Value = $game.Away
    $ws.Cells.Item($r, 2).Value = $game.AwayPts
    $ws.Cells.Item($r, 3).Value = $game.Home
    $ws.Cells.Item($r, 4).Value = $game.HomePts
    $ws.Cells.Item($r, 5).Value = $game.OT
    $ws.Cells.Item($r, 6).Value = $game.Attend
    $ws.Cells.Item($r, 7).Value = $game.Arena
    $ws.Cells.Item($r, 8).Value = $game.Win
    $ws.Cells.Item($r, 9).Value = $game.Loss
}

# Match the saved workbook's view state: selection on the new last cell and
# the window scrolled so that row 803 is the top visible row (mirrors
# topLeftCell="A803" / selection A834 in the target file).
$lastCell = $ws.Range("A834")
$lastCell.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 803
$win.ScrollColumn = 1

Write-Host "Appended $($newGames.Count) rows; sheet now spans $($ws.UsedRange.Rows.Count) rows."
